$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace rows 3-12 (account / name / balance) with the updated data set.
# Row layout is: Conta (col A), Nome (col B), Saldo (col C)

$ws.Cells.Item(3, 1).Value = "004210959"
$ws.Cells.Item(3, 2).Value = "ANA"
$ws.Cells.Item(3, 3).Value = 99821.99

$ws.Cells.Item(4, 1).Value = "000330949"
$ws.Cells.Item(4, 2).Value = "RENATO"
$ws.Cells.Item(4, 3).Value = 49000

$ws.Cells.Item(5, 1).Value = "005064129"
$ws.Cells.Item(5, 2).Value = "THIAGO"
$ws.Cells.Item(5, 3).Value = 22551.67

$ws.Cells.Item(6, 1).Value = "004870019"
$ws.Cells.Item(6, 2).Value = "MARIA"
$ws.Cells.Item(6, 3).Value = 22226.68

$ws.Cells.Item(7, 1).Value = "005170415"
$ws.Cells.Item(7, 2).Value = "MONICA"
$ws.Cells.Item(7, 3).Value = 15976.5

$ws.Cells.Item(8, 1).Value = "005624274"
$ws.Cells.Item(8, 2).Value = "CLAYTON"
$ws.Cells.Item(8, 3).Value = 15270.53

$ws.Cells.Item(9, 1).Value = "004212581"
$ws.Cells.Item(9, 2).Value = "MARIA"
$ws.Cells.Item(9, 3).Value = 13176.09

$ws.Cells.Item(10, 1).Value = "005046790"
$ws.Cells.Item(10, 2).Value = "BEATRIZ"
$ws.Cells.Item(10, 3).Value = 12875.37

$ws.Cells.Item(11, 1).Value = "004267119"
$ws.Cells.Item(11, 2).Value = "ANA"
$ws.Cells.Item(11, 3).Value = 4976.5

$ws.Cells.Item(12, 1).Value = "005683532"
$ws.Cells.Item(12, 2).Value = "SYLVERSON"
$ws.Cells.Item(12, 3).Value = 1833.46
